# Add a "Save" column (H) to the s_vals worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the other header cells (reuse G1's formatting,
# i.e. same bold/border/centered style used across row 1).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save flag values for rows 2-31, taken from the committed data.
$saveValues = @(0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,1,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
